# Weekly price-data update: insert a new daily record as row 271,
# pushing the existing rows 271-326 down to 272-327.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 271 (shifts 271..326 -> 272..327)
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new record
$ws.Range("A271").Value = 7
$ws.Range("B271").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C271").Value = "Ñuble"
$ws.Range("D271").Value = 45209
$ws.Range("E271").Value = 16
$ws.Range("F271").Value = 100112040
$ws.Range("G271").Value = "Cilantro"
$ws.Range("H271").Value = "Sin especificar"
$ws.Range("I271").Value = "Primera"
$ws.Range("J271").Value = 150
$ws.Range("K271").Value = 1500
$ws.Range("L271").Value = 1500
$ws.Range("M271").Value = 1500
$ws.Range("N271").Value = "$/atado 0,5 a 1 kilo"
$ws.Range("O271").Value = "Región de Ñuble"
$ws.Range("P271").Value = 1500
$ws.Range("Q271").Value = 1
$ws.Range("R271").Value = "Hortaliza"
